$wbx = $excel.ActiveWorkbook
$ws = $wbx.Worksheets.Item("Issue Tracking")
$bullet = [char]0x2022

# --- Row 13 (Issue #11, "Firewall rules not loading..."): testing update appended, status -> Under Test (GW) ---
$e13 = $ws.Range("E13").Value()
$ws.Range("E13").Value = $e13 + "`n" + $bullet + "01/06 Tested the last report from the customer regarding the port forwarding rule not working together with the rule of LuvitRED and I have not found the issue the customer is reporting."
$ws.Range("G13").Value = "Under Test (GW)"

# --- Row 6 (Issue #4, "For the configuration above..."): testing update appended ---
$e6 = $ws.Range("E6").Value()
$ws.Range("E6").Value = $e6 + "`n" + $bullet + "01/06 Seems to be fixed on the 2.70.0 firmware. GetWireless to confirm."
$ws.Rows.Item(6).RowHeight = 105

# --- Row 7 (Issue #5, "We first tested the Gemalto..."): testing update appended, status -> Under Test (GW) ---
$e7 = $ws.Range("E7").Value()
$ws.Range("E7").Value = $e7 + "`n" + $bullet + "01/06 Engineering build provided to GetWireless. GetWireless to test."
$ws.Range("G7").Value = "Under Test (GW)"

# --- Row 9 (Issue #7, "Release notes should be adapted..."): status -> Closed ---
$ws.Range("G9").Value = "Closed"

# --- Row 12 (Issue #10, "LuvitRED crashing..."): status -> Closed, row hidden ---
$ws.Range("G12").Value = "Closed"
$ws.Rows.Item(12).Hidden = $true

# --- View state: selection moved to F8 ---
$ws.Activate()
$ws.Range("F8").Select()
